$wb = $excel.ActiveWorkbook

# --- R1: update elapsed-duration text values, remove the now-stale row 5 ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3876:00:00"
$ws1.Range("G3").Value = "15:30:00"
$ws1.Rows(5).Delete()

# --- R2: update elapsed-duration text values ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12057:00:00"
$ws2.Range("G3").Value = "3186:42:00"
$ws2.Range("G4").Value = "424:54:00"

# --- R4: update elapsed-duration text values ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2902:48:00"
$ws4.Range("G3").Value = "130:00:00"

# --- R5: update elapsed-duration text value ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "376:48:00"

# --- R6: update elapsed-duration text value ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "17:24:00"
